# agrego CLEAS para motor tercero
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "CLEAS" column (AH), with header and values ---
$ws.Range("AH1").Value = "CLEAS"
$ws.Range("AH2").Value = "No"
$ws.Range("AH3").Value = "No"
$ws.Range("AH4").Value = "No"
$ws.Range("AH5").Value = "Sí"
$ws.Range("AH6").Value = "No"
$ws.Range("AH7").Value = "No"
$ws.Range("AH8").Value = "No"
$ws.Range("AH9").Value = "No"

# --- Update NroPoliza (F) and FechaSiniestro (H) for rows 3-5 to the new claim ---
$ws.Range("F3").Value = "04104019009"
$ws.Range("H3").Value = "20/09/2021"

$ws.Range("F4").Value = "04104019009"
$ws.Range("H4").Value = "20/09/2021"

$ws.Range("F5").Value = "04104019009"
$ws.Range("H5").Value = "20/09/2021"

# --- Anio (V) column corrected to 2020 for rows 4-9 ---
$ws.Range("V4").Value = 2020
$ws.Range("V5").Value = 2020
$ws.Range("V6").Value = 2020
$ws.Range("V7").Value = 2020
$ws.Range("V8").Value = 2020
$ws.Range("V9").Value = 2020

# --- Responsabilidad (AG) updated for row 4 ---
$ws.Range("AG4").Value = "Sí"

# --- Restore the active selection/view state as left by the author ---
$ws.Range("W8").Select()
